# Adds a new column L (2020) mirroring column K (2019) data/formatting
# for the forest-area table, and updates the sheet's dimension + the
# active selection, matching the upstream "2020 column added" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to place in the new L column, row by row (row => value).
$values = @{
    4  = 2020
    5  = 5.6
    6  = 0.8
    7  = 1.9
    8  = 0.7
    9  = 0.7
    10 = 0.9
    11 = 0.3
    12 = 0.2
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Range("K$row")
    $dstCell = $ws.Range("L$row")

    # Copy column K's number formatting/style onto the new column L cell,
    # then write the new value (keeps the existing style index, e.g. s="13").
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
    $dstCell.Value = $values[$row]
}

$excel.CutCopyMode = 0

# New selection left by the edit, per the diff.
[void]$ws.Range("N5").Select()
